$wb = $excel.ActiveWorkbook

$wsTech = $wb.Worksheets.Item("Technologies")
$wsEmp  = $wb.Worksheets.Item("Employees")
$wsRole = $wb.Worksheets.Item("Roles")

# --- Bug fix: Mrunali Desai's EMP ID was wrong (duplicated Suraksha
#     Nigade's "E0123" instead of her own "E0122"-slot id). Correcting it
#     collapses the two identical "E0123" shared strings into one.
$wsEmp.Range("B3").Value = "E0123"

# --- Bug fix: employee E0124's row had an incorrect name carried over
#     ("Pratiksha Sanam") duplicated from elsewhere; clear it so the
#     bulk-upload template row for E0124 has a blank name cell again.
$wsEmp.Range("C5").Value = ""

# LibreOffice recalculates the (auto) row height once the long name text
# is removed from C5 - mirror that slightly shorter optimal row height.
$wsEmp.Rows.Item(5).RowHeight = 15.65

# --- Restore cursor / selection bookkeeping on each sheet and make the
#     Employees tab the active one (matches the saved view state).
$null = $wsTech.Range("B9").Select()
$null = $wsRole.Range("D9").Select()

$null = $wsEmp.Activate()
$null = $wsEmp.Range("C5").Select()
